$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lunch_theroost")

# Clear the unfinished "Nutrition Label" placeholders (F2:F10, formerly "waffles")
$ws.Range("F2:F10").Value = "placeholder"

# Clear the unfinished "Allergens" placeholders (C3:C10, formerly "No known priority allergens")
$ws.Range("C3:C10").Value = "placeholder"

# Update active selection to C10
$ws.Range("C10").Select()
